# Finish all the test case about product browse
# Adds five new worksheets (005_iMacsPage .. 009_MacBooksPage) after the
# existing 004_AccessoriesPage sheet, each following the same
# TestCase / Product Name / Product Prices layout used by the other pages,
# and updates the active sheet / selection state to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the five new worksheets, in order, right after the last
#    existing sheet (004_AccessoriesPage).
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$sheetNames = @(
    "005_iMacsPage",
    "006_iPadsPage",
    "007_iPhonesPage",
    "008_iPodsPage",
    "009_MacBooksPage"
)

$newSheets = @{}
foreach ($name in $sheetNames) {
    $s = $wb.Worksheets.Add($null, $afterSheet)
    $s.Name = $name
    $newSheets[$name] = $s
    $afterSheet = $s
}

# ---------------------------------------------------------------------
# 2. Fill in the header row + data rows for each new sheet.
#    Column C ("Product Prices") is stored as text (e.g. "$150.00"),
#    matching the formatting already used on 004_AccessoriesPage.
# ---------------------------------------------------------------------
function Set-ProductRow($ws, $row, $testCase, $product, $price) {
    if ($testCase -ne $null) {
        $ws.Range("A$row").Value = $testCase
    }
    $ws.Range("B$row").Value = $product
    $priceCell = $ws.Range("C$row")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
}

# --- 005_iMacsPage ---
$ws = $newSheets["005_iMacsPage"]
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Product Name"
$ws.Range("C1").Value = "Product Prices"
Set-ProductRow $ws 2 "005-Test iMacs Page" "Magic Mouse" "`$150.00"

# --- 006_iPadsPage ---
$ws = $newSheets["006_iPadsPage"]
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Product Name"
$ws.Range("C1").Value = "Product Prices"
Set-ProductRow $ws 2 "006-Test iPads Page" "Magic Mouse" "`$150.00"
Set-ProductRow $ws 3 $null "Apple iPad 2 16GB, Wi-Fi, 9.7in – Black" "`$270.00"
Set-ProductRow $ws 4 $null "Apple iPad 6 32GB (White, 3D)" "`$680.00"

# --- 007_iPhonesPage ---
$ws = $newSheets["007_iPhonesPage"]
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Product Name"
$ws.Range("C1").Value = "Product Prices"
Set-ProductRow $ws 2 "007-Test iPhones Page" "Magic Mouse" "`$150.00"
Set-ProductRow $ws 3 $null "Apple iPhone 4S 16GB SIM-Free – Black" "`$270.00"
Set-ProductRow $ws 4 $null "Apple iPhone 4S 32GB SIM-Free – White" "`$270.00"

# --- 008_iPodsPage ---
$ws = $newSheets["008_iPodsPage"]
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Product Name"
$ws.Range("C1").Value = "Product Prices"
Set-ProductRow $ws 2 "008-Test iPods Page" "Magic Mouse" "`$150.00"
Set-ProductRow $ws 3 $null "Apple iPod touch 32GB 5th Generation – Black" "`$204.00"
Set-ProductRow $ws 4 $null "Apple iPod touch Large" "`$10.00"

# --- 009_MacBooksPage ---
$ws = $newSheets["009_MacBooksPage"]
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "Product Name"
$ws.Range("C1").Value = "Product Prices"
Set-ProductRow $ws 2 "009-Test MacBooks Page" "Magic Mouse" "`$150.00"
Set-ProductRow $ws 3 $null "Apple 13-inch MacBook Pro" "`$864.00"

# ---------------------------------------------------------------------
# 3. Restore/update the selections on each sheet and the active tab.
#    004_AccessoriesPage's old "mid-edit" selection (C6) is cleared to
#    a full-table selection now that work on it is done; the new
#    008_iPodsPage becomes the active tab, parked on A2.
# ---------------------------------------------------------------------
$acc = $wb.Worksheets.Item("004_AccessoriesPage")
$acc.Activate()
$acc.Range("A1:C7").Select()

$newSheets["005_iMacsPage"].Activate()
$newSheets["005_iMacsPage"].Range("A1:C2").Select()

$newSheets["006_iPadsPage"].Activate()
$newSheets["006_iPadsPage"].Range("A1:C4").Select()

$newSheets["007_iPhonesPage"].Activate()
$newSheets["007_iPhonesPage"].Range("A1:C4").Select()

$newSheets["009_MacBooksPage"].Activate()
$newSheets["009_MacBooksPage"].Range("B8").Select()

# 008_iPodsPage is the sheet left selected/active when the file was saved.
$newSheets["008_iPodsPage"].Activate()
$newSheets["008_iPodsPage"].Range("A2").Select()
